$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Capabilities": EndPoint column (D) changes from the pCloudy device
# endpoint to the US endpoint, for both data rows (row 2 and row 3).
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capabilities")
$wsCap.Range("D2").Value = "https://us.pcloudy.com"
$wsCap.Range("D3").Value = "https://us.pcloudy.com"

# Update the selected cell on this sheet (was D8, now D7).
$wsCap.Activate()
$wsCap.Range("D7").Select()

# ---------------------------------------------------------------------------
# Sheet "DeviceList": the Samsung Android devices used for the two test
# columns (B / C) are replaced with Apple iOS devices.
# ---------------------------------------------------------------------------
$wsDev = $wb.Worksheets.Item("DeviceList")

# Device name
$wsDev.Range("B1").Value = "APPLE_iPhone11Pro_iOS_14.4.0_6ccce"
$wsDev.Range("C1").Value = "APPLE_iPhone8_iOS_14.1.0_81551"

# Version (quote-prefixed so the text-like value keeps the same cell style)
$wsDev.Range("B2").Value = "'14.4.0"
$wsDev.Range("C2").Value = "'14.1.0"

# OperatingSystem: switch from Android to iOS (quote-prefixed, same reason)
$wsDev.Range("B3").Value = "'pCloudyIOS"
$wsDev.Range("C3").Value = "'pCloudyIOS"

# Update the selected cell on this sheet (was C4, now C18).
$wsDev.Activate()
$wsDev.Range("C18").Select()
